$d = $word.ActiveDocument

$replacements = @(
    @("37×73=2701", "85×14=1190"),
    @("74×44=3256", "56×87=4872"),
    @("60×28=1680", "68×36=2448"),
    @("26×13=338", "31×41=1271"),
    @("32×77=2464", "60×36=2160"),
    @("50×72=3600", "16×16=256"),
    @("41×86=3526", "89×66=5874"),
    @("53×25=1325", "43×43=1849"),
    @("98×79=7742", "18×15=270"),
    @("19×80=1520", "79×22=1738"),
    @("96×81=7776", "82×51=4182"),
    @("20×37=740", "92×77=7084"),
    @("76×42=3192", "63×73=4599"),
    @("67×25=1675", "29×87=2523"),
    @("88×17=1496", "38×13=494"),
    @("72×54=3888", "17×57=969"),
    @("95×63=5985", "22×30=660"),
    @("14×59=826", "56×88=4928"),
    @("28×65=1820", "68×70=4760"),
    @("81×92=7452", "24×28=672"),
    @("36×50=1800", "62×81=5022"),
    @("26×67=1742", "99×40=3960"),
    @("93×98=9114", "93×84=7812"),
    @("13×83=1079", "66×80=5280"),
    @("69×13=897", "27×57=1539")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
